$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "MIE(5.3058492471983385, 4.297012431873066, -6.549959678498338, 9.521250737270659)"
$ws.Range("C2").Value = "NIG(0.7124466850086493, 0.46998261716199063, 6.015706881425392, 4.8627006721890975)"
$ws.Range("D2").Value = "NCT(2.950481764569178, 1.5684429365770551, -0.8381683677234233, 2.4222175445306253)"
$ws.Range("E2").Value = "NIG(1.0455584505516207, 0.8055202023944532, 4.800455262518598, 5.050767645649259)"
